# Add two new question rows (14 & 15) to the "Question List" worksheet,
# mirroring the existing rows in the table (module / description / asker / date),
# then update the row heights for the new wrapped text and move the active
# selection to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: new question about application management ---
$ws.Range("B14").Value = "应用管理"
$ws.Range("C14").Value = "没有找到任何有关应用管理接口的文档。请提供"
# Copy the asker / date cells from row 13 so the values are stored as plain
# text (matching the existing "耿晓红" / "2015.11.15" text already used in
# this sheet) instead of being auto-converted into a date serial number.
$ws.Range("D13").Copy($ws.Range("D14"))
$ws.Range("E13").Copy($ws.Range("E14"))

# --- Row 15: follow-up question about application management ---
$ws.Range("B15").Value = "应用管理"
$ws.Range("C15").Value = "应用管理的菜单在页面上只有专题管理和精品课程管理，我们只需要这两个吗？"
$ws.Range("D13").Copy($ws.Range("D15"))
$ws.Range("E13").Copy($ws.Range("E15"))

# The new text wraps onto multiple lines, so the rows need to grow taller.
$ws.Rows.Item(14).RowHeight = 27
$ws.Rows.Item(15).RowHeight = 40.5

# Move the selection to the newly-added second question row.
$ws.Activate()
$ws.Range("D15:E15").Select()
